# ---------------------------------------------------------------------------
# Revise antibody template 1.2
#
# Applies the "Revise antibody template 1.2" edit to the CoVIC-DB antibodies
# submission workbook:
#   - Instructions sheet: rewritten with a version line + a column
#     reference table (Antibody name / Host / Isotype / Light chain /
#     Heavy chain germline / Antibody details / Structural data /
#     Antibody comment) in columns A/B.
#   - Antibodies sheet: adds "Light chain", "Heavy chain germline",
#     "Structural data" columns (and reorders "Antibody details" /
#     "Antibody comment"), plus matching list validations.
#   - Terminology sheet: adds "Light chain" and "Heavy chain germline"
#     lookup columns.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Excel's ColumnWidth property is expressed in "characters of the Normal
# style font", which the file format then stores as character-width + a
# fixed ~5px padding (0.8333... chars at the default font). Centralise the
# conversion so the widths we request match the widths the diff expects.
function Set-ColWidth($ws, $colIndex, $targetWidth) {
    $ws.Columns.Item($colIndex).ColumnWidth = $targetWidth - (5.0 / 6.0)
}

# Copies formatting only (no value) from $srcRange onto $dstRange, cell by
# cell, via Copy + PasteSpecial(xlPasteFormats). Used both to (re)apply the
# bold header style (reusing the sheet's existing bold style instead of
# minting a new one) and to force an "empty but present" cell to survive
# the save round-trip.
function Copy-Format($srcRange, $dstRange) {
    $srcRange.Copy() | Out-Null
    $dstRange.PasteSpecial(-4122) | Out-Null
}

# ===========================================================================
# Sheet 1: Instructions
# ===========================================================================
$ws1 = $wb.Worksheets.Item("Instructions")
$ws1.Unprotect()

$ws1.Cells.Clear() | Out-Null

# A1 is the bold title; keep it bold and use it as the "bold style" source
# for every other header cell below, and as the "default style" source
# is taken from a freshly-cleared cell (A1 before we restyle it).
$defaultStyleSrc = $ws1.Range("Z1")

$ws1.Range("A1").Value = "CoVIC-DB Antibodies Submission"
$ws1.Range("A1").Font.Bold = $true
$boldSrc = $ws1.Range("A1")

# B1: empty placeholder cell (present, no value) next to the title.
Copy-Format $defaultStyleSrc $ws1.Range("B1")

$ws1.Range("A2").Value = "Version 1.2"

# A3: empty row.
Copy-Format $defaultStyleSrc $ws1.Range("A3")

$ws1.Range("A4").Value = "Add your antibodies to the 'Antibodies' sheet. Do not edit the other sheets."

$headerRows = @(6, 7, 8, 9, 10, 11, 19, 22)
$labels = @{
    6  = "Antibody name"
    7  = "Host"
    8  = "Isotype"
    9  = "Light chain"
    10 = "Heavy chain germline"
    11 = "Antibody details"
    19 = "Structural data"
    22 = "Antibody comment"
}
$descriptions = @{
    6  = "Your preferred code name for the antibody"
    7  = "Specify the host species that is the source of the antibody"
    8  = "Specify the antibody isotype, if known"
    9  = "Specify the antibody light chain, if known (kappa or lambda)"
    10 = "Specify the antibody heavy chain germline gene, if known"
    11 = "Measurements or characteristics of the antibody."
    19 = "Would you like structural analyses of this antibody?"
    22 = "Please provide any other details about the antibody."
}

foreach ($r in $headerRows) {
    $ws1.Range("A$r").Value = $labels[$r]
    Copy-Format $boldSrc $ws1.Range("A$r")
    $ws1.Range("B$r").Value = $descriptions[$r]
}

# Rows that only carry explanatory text in column B (column A left blank).
$bOnly = @{
    12 = "This column is optional, and meant to capture data you might have on the antibody."
    13 = "These data will not be released to the partner reference labs that will perform the analyses."
    14 = "For example:"
    15 = "- Affinity: Spike protein binding affinity; inhibition of ACE2 binding; ELISA for Spike "
    16 = "- Neutralization: IC50 value"
    17 = "- Neutralization assay platform"
    18 = "- Epitope: Binning or competition data"
    20 = "If no, leave blank."
    21 = "If yes, rank the antibodies in order of priority, starting with '1' for the highest priority."
}
foreach ($r in $bOnly.Keys) {
    Copy-Format $defaultStyleSrc $ws1.Range("A$r")
    $ws1.Range("B$r").Value = $bOnly[$r]
}

Set-ColWidth $ws1 1 18
Set-ColWidth $ws1 2 70

$ws1.Protect() | Out-Null

# ===========================================================================
# Sheet 2: Antibodies
# ===========================================================================
$ws2 = $wb.Worksheets.Item("Antibodies")

$ws2.Range("D1").Value = "Light chain"
$ws2.Range("E1").Value = "Heavy chain germline"
$ws2.Range("F1").Value = "Antibody details"
$ws2.Range("G1").Value = "Structural data"
$ws2.Range("H1").Value = "Antibody comment"

foreach ($addr in @("D1", "E1", "F1", "G1", "H1")) {
    Copy-Format $ws2.Range("A1") $ws2.Range($addr)
}

# Keep row 2 present (matches the original placeholder <row r="2"/>).
$ws1.Range("A3").Copy() | Out-Null
$ws2.Range("A2").PasteSpecial(-4122) | Out-Null

Set-ColWidth $ws2 1 15
Set-ColWidth $ws2 2 15
Set-ColWidth $ws2 3 15
Set-ColWidth $ws2 4 15
Set-ColWidth $ws2 5 20
Set-ColWidth $ws2 6 16
Set-ColWidth $ws2 7 15
Set-ColWidth $ws2 8 16

$ws2.Range("D2:D100").Validation.Add(3, 1, 1, "=Terminology!`$C`$2:`$C`$4") | Out-Null
$ws2.Range("E2:E100").Validation.Add(3, 1, 1, "=Terminology!`$D`$2:`$D`$12") | Out-Null

# ===========================================================================
# Sheet 3: Terminology
# ===========================================================================
$ws3 = $wb.Worksheets.Item("Terminology")
$ws3.Unprotect()

$ws3.Range("C1").Value = "Light chain"
$ws3.Range("D1").Value = "Heavy chain germline"
foreach ($addr in @("C1", "D1")) {
    Copy-Format $ws3.Range("A1") $ws3.Range($addr)
}

$lightChain = @{
    2 = "kappa"
    3 = "lambda"
    4 = "unknown"
}
foreach ($r in $lightChain.Keys) {
    $ws3.Range("C$r").Value = $lightChain[$r]
}

$germline = @{
    2  = "IGHV1-8"
    3  = "IGHV1-18"
    4  = "IGHV2-5"
    5  = "IGHV3-7"
    6  = "IGHV3-11"
    7  = "IGHV3-21"
    8  = "IGHV3-23"
    9  = "IGHV4-39"
    10 = "IGHV4-59"
    11 = "IGHV5-51"
    12 = "IGHV6-1"
}
foreach ($r in $germline.Keys) {
    $ws3.Range("D$r").Value = $germline[$r]
}

# C5:C15 and D13:D15 stay blank, but present (matches the diff's empty
# inlineStr placeholder cells). Copy formatting (not value) from an
# untouched, default-styled cell so these stay unstyled, like their
# neighbours.
$ws3DefaultStyleSrc = $ws3.Range("Z1")
foreach ($r in 5..15) {
    Copy-Format $ws3DefaultStyleSrc $ws3.Range("C$r")
}
foreach ($r in 13..15) {
    Copy-Format $ws3DefaultStyleSrc $ws3.Range("D$r")
}

Set-ColWidth $ws3 1 15
Set-ColWidth $ws3 2 15
Set-ColWidth $ws3 3 15
Set-ColWidth $ws3 4 20

$ws3.Protect() | Out-Null

Write-Host "Revised antibody template to 1.2"
